# Auto-generated script applying scheduled-runner price/profit refresh
# across the Excalibur_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 883.3333
$ws.Range("I12").Value = 740.5
$ws.Range("J12").Value = 954.75
$ws.Range("K12").Value = 740.5
$ws.Range("L12").Value = 954.75
$ws.Range("M12").Value = -570.5
$ws.Range("N12").Value = -1294.75
$ws.Range("H17").Value = 1001.5857
$ws.Range("J17").Value = 1011.8788
$ws.Range("L17").Value = 3035.6364
$ws.Range("N17").Value = -3371.6364
$ws.Range("H32").Value = 3757.75
$ws.Range("I32").Value = 3479.4285
$ws.Range("K32").Value = 3479.4285
$ws.Range("M32").Value = -3153.4285
$ws.Range("H33").Value = 253
$ws.Range("I33").Value = 249.08333
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 249.08333
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -20.08332999999999
$ws.Range("N33").Value = -758
$ws.Range("H55").Value = 469.8
$ws.Range("J55").Value = 526
$ws.Range("L55").Value = 526
$ws.Range("N55").Value = -954
$ws.Range("H64").Value = 7411.75
$ws.Range("I64").Value = 4990.6
$ws.Range("J64").Value = 9141.143
$ws.Range("K64").Value = 4990.6
$ws.Range("L64").Value = 9141.143
$ws.Range("M64").Value = -4742.6
$ws.Range("N64").Value = -9637.143
$ws.Range("H67").Value = 7411.75
$ws.Range("I67").Value = 4990.6
$ws.Range("J67").Value = 9141.143
$ws.Range("K67").Value = 4990.6
$ws.Range("L67").Value = 9141.143
$ws.Range("M67").Value = -4132.6
$ws.Range("N67").Value = -10857.143
$ws.Range("H86").Value = 2030.05
$ws.Range("I86").Value = 1014.5714
$ws.Range("J86").Value = 2576.8462
$ws.Range("K86").Value = 1014.5714
$ws.Range("L86").Value = 2576.8462
$ws.Range("M86").Value = 108.4286
$ws.Range("N86").Value = -4822.8462
$ws.Range("H89").Value = 2030.05
$ws.Range("I89").Value = 1014.5714
$ws.Range("J89").Value = 2576.8462
$ws.Range("K89").Value = 5072.857
$ws.Range("L89").Value = 12884.231
$ws.Range("M89").Value = 543.143
$ws.Range("N89").Value = -24116.231
$ws.Range("H112").Value = 1118.6757
$ws.Range("J112").Value = 1135.0588
$ws.Range("L112").Value = 3405.1764
$ws.Range("N112").Value = -5621.1764
$ws.Range("H113").Value = 3483.1667
$ws.Range("I113").Value = 3224.75
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3224.75
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 29.25
$ws.Range("N113").Value = -10508
$ws.Range("H116").Value = 214763.2
$ws.Range("I116").Value = 99005.5
$ws.Range("K116").Value = 99005.5
$ws.Range("M116").Value = -95563.5
$ws.Range("H137").Value = 956084.0600000001
$ws.Range("I137").Value = 950.53845
$ws.Range("J137").Value = 1547357.2
$ws.Range("K137").Value = 2851.61535
$ws.Range("L137").Value = 4642071.6
$ws.Range("M137").Value = -301.61535
$ws.Range("N137").Value = -4647171.6
$ws.Range("H138").Value = 6850.24
$ws.Range("I138").Value = 13435.474
$ws.Range("J138").Value = 2814.1292
$ws.Range("K138").Value = 40306.422
$ws.Range("L138").Value = 8442.3876
$ws.Range("M138").Value = -35166.422
$ws.Range("N138").Value = -18722.3876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7983.7
$ws.Range("I2").Value = 6666
$ws.Range("K2").Value = 6666
$ws.Range("M2").Value = -6553
$ws.Range("H61").Value = 856410.6
$ws.Range("I61").Value = 1112067.4
$ws.Range("K61").Value = 1112067.4
$ws.Range("M61").Value = -1111855.4
$ws.Range("H116").Value = 7983.7
$ws.Range("I116").Value = 6666
$ws.Range("K116").Value = 6666
$ws.Range("M116").Value = -4372
$ws.Range("H132").Value = 373169.44
$ws.Range("I132").Value = 449443.6
$ws.Range("K132").Value = 1348330.8
$ws.Range("M132").Value = -1345800.8
$ws.Range("H136").Value = 856410.6
$ws.Range("I136").Value = 1112067.4
$ws.Range("K136").Value = 3336202.2
$ws.Range("M136").Value = -3333652.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7983.7
$ws.Range("I3").Value = 6666
$ws.Range("K3").Value = 6666
$ws.Range("M3").Value = -6552
$ws.Range("H20").Value = 1475.7646
$ws.Range("I20").Value = 1196.3636
$ws.Range("K20").Value = 1196.3636
$ws.Range("M20").Value = -949.3635999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 100
$ws.Range("K3").Value = 100
$ws.Range("M3").Value = 13
$ws.Range("H31").Value = 132369.62
$ws.Range("I31").Value = 185244.17
$ws.Range("K31").Value = 185244.17
$ws.Range("M31").Value = -184949.17
$ws.Range("H34").Value = 132369.62
$ws.Range("I34").Value = 185244.17
$ws.Range("K34").Value = 185244.17
$ws.Range("M34").Value = -185042.17
$ws.Range("H86").Value = 140662.06
$ws.Range("I86").Value = 7303.5
$ws.Range("K86").Value = 7303.5
$ws.Range("M86").Value = -6180.5
$ws.Range("H89").Value = 140662.06
$ws.Range("I89").Value = 7303.5
$ws.Range("K89").Value = 36517.5
$ws.Range("M89").Value = -30901.5
$ws.Range("H132").Value = 23595790
$ws.Range("I132").Value = 28584974
$ws.Range("J132").Value = 13894599
$ws.Range("K132").Value = 85754922
$ws.Range("L132").Value = 41683797
$ws.Range("M132").Value = -85752392
$ws.Range("N132").Value = -41688857
$ws.Range("H134").Value = 11463.849
$ws.Range("I134").Value = 12510.173
$ws.Range("K134").Value = 37530.519
$ws.Range("M134").Value = -34995.519
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1160.5
$ws.Range("I129").Value = 795.2727
$ws.Range("J129").Value = 2499.6667
$ws.Range("K129").Value = 2385.8181
$ws.Range("L129").Value = 7499.000100000001
$ws.Range("M129").Value = 2614.1819
$ws.Range("N129").Value = -17499.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 235.375
$ws.Range("I2").Value = 120.55556
$ws.Range("J2").Value = 383
$ws.Range("K2").Value = 120.55556
$ws.Range("L2").Value = 383
$ws.Range("M2").Value = -7.55556
$ws.Range("N2").Value = -609
$ws.Range("H17").Value = 775.75
$ws.Range("I17").Value = 743.7143
$ws.Range("K17").Value = 743.7143
$ws.Range("M17").Value = -575.7143
$ws.Range("H75").Value = 53249.25
$ws.Range("I75").Value = 40000
$ws.Range("J75").Value = 57665.668
$ws.Range("K75").Value = 40000
$ws.Range("L75").Value = 57665.668
$ws.Range("M75").Value = -39126
$ws.Range("N75").Value = -59413.668
$ws.Range("H78").Value = 53249.25
$ws.Range("I78").Value = 40000
$ws.Range("J78").Value = 57665.668
$ws.Range("K78").Value = 120000
$ws.Range("L78").Value = 172997.004
$ws.Range("M78").Value = -115632
$ws.Range("N78").Value = -181733.004
$ws.Range("H132").Value = 11776346
$ws.Range("I132").Value = 16072564
$ws.Range("J132").Value = 8444.348
$ws.Range("K132").Value = 48217692
$ws.Range("L132").Value = 25333.044
$ws.Range("M132").Value = -48215162
$ws.Range("N132").Value = -30393.044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 446.58334
$ws.Range("I16").Value = 431.92856
$ws.Range("J16").Value = 467.1
$ws.Range("K16").Value = 431.92856
$ws.Range("L16").Value = 467.1
$ws.Range("M16").Value = -261.92856
$ws.Range("N16").Value = -807.1
$ws.Range("H22").Value = 436.9375
$ws.Range("J22").Value = 449.625
$ws.Range("L22").Value = 449.625
$ws.Range("N22").Value = -1039.625
$ws.Range("H27").Value = 436.9375
$ws.Range("J27").Value = 449.625
$ws.Range("L27").Value = 449.625
$ws.Range("N27").Value = -663.625
$ws.Range("H68").Value = 2760
$ws.Range("J68").Value = 2766.6667
$ws.Range("L68").Value = 2766.6667
$ws.Range("N68").Value = -4264.6667
$ws.Range("H71").Value = 2760
$ws.Range("J71").Value = 2766.6667
$ws.Range("L71").Value = 13833.3335
$ws.Range("N71").Value = -21321.3335
$ws.Range("H104").Value = 70547.8
$ws.Range("J104").Value = 70547.8
$ws.Range("L104").Value = 70547.8
$ws.Range("N104").Value = -77535.8
$ws.Range("H122").Value = 5600
$ws.Range("I122").Value = 5090.909
$ws.Range("K122").Value = 15272.727
$ws.Range("M122").Value = -12822.727
$ws.Range("H132").Value = 1090117.5
$ws.Range("I132").Value = 1394009.8
$ws.Range("K132").Value = 4182029.4
$ws.Range("M132").Value = -4179499.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6294
$ws.Range("I81").Value = 6008.9
$ws.Range("J81").Value = 6610.778
$ws.Range("K81").Value = 12017.8
$ws.Range("L81").Value = 13221.556
$ws.Range("M81").Value = -10956.8
$ws.Range("N81").Value = -15343.556
$ws.Range("H84").Value = 6294
$ws.Range("I84").Value = 6008.9
$ws.Range("J84").Value = 6610.778
$ws.Range("K84").Value = 60089
$ws.Range("L84").Value = 66107.78
$ws.Range("M84").Value = -54785
$ws.Range("N84").Value = -76715.78
$ws.Range("H122").Value = 3344.9285
$ws.Range("I122").Value = 3160.1155
$ws.Range("J122").Value = 5747.5
$ws.Range("K122").Value = 9480.3465
$ws.Range("L122").Value = 17242.5
$ws.Range("M122").Value = -7030.3465
$ws.Range("H132").Value = 4377670.5
$ws.Range("I132").Value = 5751480.5
$ws.Range("K132").Value = 17254441.5
$ws.Range("M132").Value = -17251911.5
$ws.Range("H136").Value = 1348025.5
$ws.Range("I136").Value = 1733241
$ws.Range("J136").Value = 63973.668
$ws.Range("K136").Value = 5199723
$ws.Range("L136").Value = 191921.004
$ws.Range("M136").Value = -5197173
$ws.Range("N136").Value = -197021.004

